# Apply cryptos.xlsx update (values & percentage changes) via Excel COM interop
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text / percentage cell updates (safe as text, no numeric coercion) ---
$ws.Range("D2").Value = "36.613.49"
$ws.Range("E2").Value = "  -0.08%  "
$ws.Range("D3").Value = "1.961.19"
$ws.Range("E3").Value = "  +0.03%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("E6").Value = "  -0.55%  "
$ws.Range("E7").Value = "  -0.43%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("E9").Value = "  +2.56%  "
$ws.Range("E10").Value = "  -5.90%  "
$ws.Range("E11").Value = "  -0.85%  "
$ws.Range("E12").Value = "  -2.61%  "
$ws.Range("E13").Value = "  -0.19%  "
$ws.Range("D14").Value = "2.246.54"
$ws.Range("E14").Value = "  -0.21%  "
$ws.Range("E15").Value = "  -0.28%  "
$ws.Range("E16").Value = "  +0.40%  "
$ws.Range("D17").Value = "1.961.37"
$ws.Range("E17").Value = "  +0.12%  "
$ws.Range("D18").Value = "36.517.06"
$ws.Range("E18").Value = "  -0.14%  "
$ws.Range("E19").Value = "  -0.63%  "
$ws.Range("D20").Value = "0.0₃0856"
$ws.Range("E20").Value = "  -2.41%  "
$ws.Range("E21").Value = "  -0.90%  "
$ws.Range("E22").Value = "  -0.84%  "
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("E24").Value = "  -1.59%  "
$ws.Range("E25").Value = "  +1.73%  "
$ws.Range("E26").Value = "  -1.91%  "
$ws.Range("E27").Value = "  +1.52%  "
$ws.Range("E28").Value = "  -1.32%  "
$ws.Range("E29").Value = "  -0.92%  "
$ws.Range("E30").Value = "  +1.16%  "
$ws.Range("E31").Value = "  -3.21%  "
$ws.Range("E32").Value = "  -0.78%  "
$ws.Range("E33").Value = "  -3.63%  "
$ws.Range("E34").Value = "  -0.11%  "
$ws.Range("E35").Value = "  -0.28%  "
$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("E36").Value = "  +1.76%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("E37").Value = "  +10.15%  "
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("E39").Value = "  -11.88%  "
$ws.Range("E40").Value = "  -1.68%  "
$ws.Range("E41").Value = "  +0.96%  "
$ws.Range("E42").Value = "  -1.62%  "
$ws.Range("E43").Value = "  +0.17%  "
$ws.Range("E44").Value = "  -2.67%  "
$ws.Range("D45").Value = "1.368.43"
$ws.Range("E45").Value = "  +0.73%  "
$ws.Range("E46").Value = "  -1.37%  "
$ws.Range("E47").Value = "  -0.88%  "
$ws.Range("E48").Value = "  -1.73%  "
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("D50").Value = "2.138.06"
$ws.Range("E50").Value = "  -0.18%  "
$ws.Range("E51").Value = "  -5.14%  "

# --- Price cells whose new text looks like a plain number (e.g. "0.618") ---
# Excel auto-converts such strings to numeric values, which would change the
# cell's stored type/representation. Force the cell to Text format first so the
# literal string is preserved exactly, then restore the default "Normal" style so
# no extra formatting/style is left behind on the cell.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.618"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.82"
$ws.Range("D7").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0808"
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.19"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.832"
$ws.Range("D13").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "13.75"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.30"
$ws.Range("D16").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.81"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "228.67"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.06"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.46"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.27"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.139"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "160.59"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.47"
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0620"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.33"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.24"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.39"
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.71"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0984"
$ws.Range("D40").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "87.98"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.16"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.83"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "43.77"
$ws.Range("D51").Style = "Normal"
